$d = $word.ActiveDocument
$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main' xmlns:w14='http://schemas.microsoft.com/office/word/2010/wordml'"

function Merge-Range($range, [string]$finalText) {
    # Forces a real run-merge (identical-text replace is a no-op in this
    # engine) by first writing a distinct placeholder, then restoring the
    # exact original text; Word collapses the touched runs into one.
    $start = $range.Start
    $orig = $range.Text
    $range.Text = $orig + "~~TEMP~~"
    $r2 = $d.Range($start, $start + $orig.Length + 8)
    $r2.Text = $finalText
}

# --- 1. "Sprint BackLog #3": merge the " #" and "3" runs -> " #3" ---
$p1 = $d.Paragraphs(1)
$f1 = $p1.Range.Duplicate
$f1.Find.Execute(" #3", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
Merge-Range $f1 " #3"

# --- 2. Date line: merge the "17" and "/01/" runs -> "17/01/" ---
$p2 = $d.Paragraphs(2)
$f2 = $p2.Range.Duplicate
$f2.Find.Execute("17/01/", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
Merge-Range $f2 "17/01/"

# --- 3. Bold "Sviluppo di tutti i requisiti funzionali:" (incl. paragraph mark) ---
$p5 = $d.Paragraphs(5)
$xml5 = "<w:p $wNs w14:paraId='15C7FDB7' w14:textId='7D68FF96' w:rsidR='007F2FD0' w:rsidRPr='007F2FD0' w:rsidRDefault='007F2FD0' w:rsidP='007F2FD0'><w:pPr><w:pStyle w:val='Paragrafoelenco'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='2'/></w:numPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Sviluppo di tutti i requisiti funzionali</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>:</w:t></w:r></w:p>"
$p5.Range.InsertXML($xml5)

# --- 4. Merge + bold "Sviluppo di tutti i requisiti non funzionali e tecnologici:" ---
$p16 = $d.Paragraphs(16)
$xml16 = "<w:p $wNs w14:paraId='0D71C77D' w14:textId='16B6A336' w:rsidR='007F2FD0' w:rsidRPr='007F2FD0' w:rsidRDefault='007F2FD0' w:rsidP='007F2FD0'><w:pPr><w:pStyle w:val='Paragrafoelenco'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='2'/></w:numPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Sviluppo di tutti i requisiti non funzionali e tecnologici:</w:t></w:r></w:p>"
$p16.Range.InsertXML($xml16)

# --- 5. Merge "l" + "'ATM deve essere..." runs into one ---
$p18 = $d.Paragraphs(18)
$f18 = $p18.Range.Duplicate
$f18.Find.Execute("l’ATM deve essere il più veloce possibile nella risposta e nell’interazione con il cliente", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
Merge-Range $f18 "l’ATM deve essere il più veloce possibile nella risposta e nell’interazione con il cliente"

# --- 6. Remove the "Scadenze" block: empty para, "Scadenze " para, tab para, empty-tabs para ---
function Find-ParaIndexByText([string]$needle) {
    $i = 0
    foreach ($p in $d.Paragraphs) {
        $i = $i + 1
        if ($p.Range.Text -like "*${needle}*") {
            return $i
        }
    }
    return -1
}

$scadIdx = Find-ParaIndexByText "Scadenze"
$pPrev = $d.Paragraphs($scadIdx - 1)   # blank paragraph right before "Scadenze "
$pScad = $d.Paragraphs($scadIdx)       # "Scadenze "
$pTab = $d.Paragraphs($scadIdx + 1)    # paragraph holding just a tab
$pEmptyTabs = $d.Paragraphs($scadIdx + 2) # empty paragraph still carrying the tab stop def

$delRange = $d.Range($pPrev.Range.Start, $pEmptyTabs.Range.End)
$delRange.Delete()

Write-Output "Done"
